# "Finish how old is Teddy" — update the JS101 weekly log: row 60 went
# from 0.5 hours / "Finished 1 small problem" to 1 hour / "Finished 2
# small problems" (the now-unused "Finished 1 small problem" shared
# string disappears on save since nothing references it anymore).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 60: Hours 0.5 -> 1, Notes "Finished 1 small problem" -> "Finished 2 small problems"
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = "Finished 2 small problems"

# Move the active selection from D60 to D61, matching the saved view state.
$ws.Range("D61").Select()

$wb.Save()
